{"js": "// Replace the 25 \"NN\u00f7N=\" division-drill answers in the worksheet table.\n// Each entry is [oldText, newText]; entries are applied in the same order\n// they appear in the document so that a value which is both a target\n// (\"before\") for one cell and a replacement (\"after\") for another\n// (e.g. \"99\u00f79=\") never gets double-replaced.\nconst replacements = [\n  [\"92\u00f74=\", \"69\u00f72=\"],\n  [\"53\u00f76=\", \"43\u00f76=\"],\n  [\"54\u00f77=\", \"53\u00f78=\"],\n  [\"21\u00f75=\", \"49\u00f72=\"],\n  [\"12\u00f74=\", \"38\u00f72=\"],\n  [\"72\u00f79=\", \"23\u00f76=\"],\n  [\"99\u00f79=\", \"51\u00f74=\"],\n  [\"82\u00f74=\", \"62\u00f74=\"],\n  [\"70\u00f75=\", \"22\u00f75=\"],\n  [\"57\u00f73=\", \"15\u00f76=\"],\n  [\"27\u00f76=\", \"92\u00f75=\"],\n  [\"28\u00f74=\", \"91\u00f72=\"],\n  [\"78\u00f74=\", \"81\u00f79=\"],\n  [\"24\u00f72=\", \"62\u00f78=\"],\n  [\"59\u00f72=\", \"99\u00f79=\"],\n  [\"89\u00f75=\", \"72\u00f74=\"],\n  [\"89\u00f73=\", \"10\u00f75=\"],\n  [\"87\u00f79=\", \"42\u00f78=\"],\n  [\"20\u00f72=\", \"26\u00f73=\"],\n  [\"30\u00f78=\", \"75\u00f77=\"],\n  [\"15\u00f74=\", \"50\u00f74=\"],\n  [\"68\u00f73=\", \"99\u00f76=\"],\n  [\"64\u00f74=\", \"75\u00f75=\"],\n  [\"77\u00f72=\", \"30\u00f75=\"],\n  [\"28\u00f75=\", \"94\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NN\u00f7N=\" division-drill answers in the worksheet table.\n# Each pair is (oldText, newText); applied in the same order they appear in\n# the document so that a value which is both a target (\"before\") for one\n# cell and a replacement (\"after\") for another (e.g. \"99\u00f79=\") never gets\n# double-replaced.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"92\u00f74=\", \"69\u00f72=\"),\n  @(\"53\u00f76=\", \"43\u00f76=\"),\n  @(\"54\u00f77=\", \"53\u00f78=\"),\n  @(\"21\u00f75=\", \"49\u00f72=\"),\n  @(\"12\u00f74=\", \"38\u00f72=\"),\n  @(\"72\u00f79=\", \"23\u00f76=\"),\n  @(\"99\u00f79=\", \"51\u00f74=\"),\n  @(\"82\u00f74=\", \"62\u00f74=\"),\n  @(\"70\u00f75=\", \"22\u00f75=\"),\n  @(\"57\u00f73=\", \"15\u00f76=\"),\n  @(\"27\u00f76=\", \"92\u00f75=\"),\n  @(\"28\u00f74=\", \"91\u00f72=\"),\n  @(\"78\u00f74=\", \"81\u00f79=\"),\n  @(\"24\u00f72=\", \"62\u00f78=\"),\n  @(\"59\u00f72=\", \"99\u00f79=\"),\n  @(\"89\u00f75=\", \"72\u00f74=\"),\n  @(\"89\u00f73=\", \"10\u00f75=\"),\n  @(\"87\u00f79=\", \"42\u00f78=\"),\n  @(\"20\u00f72=\", \"26\u00f73=\"),\n  @(\"30\u00f78=\", \"75\u00f77=\"),\n  @(\"15\u00f74=\", \"50\u00f74=\"),\n  @(\"68\u00f73=\", \"99\u00f76=\"),\n  @(\"64\u00f74=\", \"75\u00f75=\"),\n  @(\"77\u00f72=\", \"30\u00f75=\"),\n  @(\"28\u00f75=\", \"94\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  # MatchAllWordForms, Forward, Wrap(0=wdFindStop), Format, ReplaceWith,\n  # Replace(2=wdReplaceAll)\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 0, $false, $null, 2) | Out-Null\n}\n"}
